$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel;
# force Text format first so the literal string (with its trailing zeros /
# exact decimal digits) is preserved, matching the source data feed.
$textCells = @(
    'D5',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D14',
    'D15',
    'D16',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D25',
    'D26',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '30.478.61'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '1.911.28'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '239.66'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.4782'
$ws.Range('E7').Value = '  -2.54%  '
$ws.Range('D8').Value = '0.2839'
$ws.Range('E8').Value = '  -3.59%  '
$ws.Range('D9').Value = '0.06689'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').Value = '19.38'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').Value = '102.44'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('D12').Value = '0.07785'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '1.927.49'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').Value = '5.200'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').Value = '0.6692'
$ws.Range('E15').Value = '  -4.51%  '
$ws.Range('D16').Value = '279.36'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('D17').Value = '30.530.68'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').Value = '0.000007469'
$ws.Range('E19').Value = '  -3.41%  '
$ws.Range('D20').Value = '12.62'
$ws.Range('E20').Value = '  -3.62%  '
$ws.Range('D21').Value = '5.388'
$ws.Range('E21').Value = '  -3.67%  '
$ws.Range('D22').Value = '0.4686'
$ws.Range('E22').Value = '  -6.64%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '6.293'
$ws.Range('E24').Value = '  -3.73%  '
$ws.Range('D25').Value = '9.347'
$ws.Range('E25').Value = '  -5.35%  '
$ws.Range('D26').Value = '167.33'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').Value = '2.077'
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('D29').Value = '1.383'
$ws.Range('E29').Value = '  -0.67%  '
$ws.Range('D30').Value = '0.09954'
$ws.Range('E30').Value = '  -4.79%  '
$ws.Range('D31').Value = '4.579'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').Value = '1.514'
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').Value = '0.04719'
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('D35').Value = '0.7241'
$ws.Range('E35').Value = '  -4.99%  '
$ws.Range('D36').Value = '1.109'
$ws.Range('E36').Value = '  -3.75%  '
$ws.Range('D37').Value = '2.712'
$ws.Range('E37').Value = '  -0.87%  '
$ws.Range('D38').Value = '0.01900'
$ws.Range('E38').Value = '  -5.59%  '
$ws.Range('D39').Value = '2.617'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').Value = '6.317'
$ws.Range('E40').Value = '  -3.15%  '
$ws.Range('D41').Value = '73.87'
$ws.Range('E41').Value = '  -5.66%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8623'
$ws.Range('E42').Value = '  -5.36%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.954'
$ws.Range('E43').Value = '  -7.02%  '
$ws.Range('D44').Value = '105.89'
$ws.Range('E44').Value = '  -2.09%  '
$ws.Range('D45').Value = '0.4245'
$ws.Range('E45').Value = '  -4.53%  '
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').Value = '7.379'
$ws.Range('E47').Value = '  -4.40%  '
$ws.Range('D48').Value = '963.13'
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('D49').Value = '0.1200'
$ws.Range('E49').Value = '  -4.14%  '
$ws.Range('D50').Value = '34.57'
$ws.Range('E50').Value = '  -4.75%  '
$ws.Range('D51').Value = '0.05797'
$ws.Range('E51').Value = '  +0.46%  '
